# "10 Apr 2025 Update"
# - Reworks the task list on the "2. Dev" sheet: several rows were
#   renumbered/shifted, a new "3.2.4 Add button to export" task was
#   inserted, and several new multilateral/chart related tasks were
#   appended at the bottom (3.23.4, 3.23.5, TVD vs MD plot, 3.25).
# - Moves the active/selected tab from "3. Well Builder" to "2. Dev".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2. Dev")
$wellBuilder = $wb.Worksheets.Item("3. Well Builder")

# --- Clear out the old task rows (8-31) so we can rebuild them cleanly ---
$ws.Range("A8:N31").ClearContents()
$ws.Range("A8:N31").ClearFormats()

# --- Row 8 / 9 : new "3.2.4 Add button to export" task, checkbox list shifts up ---
$ws.Range('C8').Value = '3.2.4'
$ws.Range('D8').Value = 'Add button to export'
$ws.Range('N8').Value = 'Show well check boxes'
$ws.Range('N9').Value = 'Highlight well check boxes'

# --- Rows 10-15 : existing tasks shift up by one slot ---
$ws.Range('A10').Value = 'X'
$ws.Range('B10').Value = 3.4
$ws.Range('C10').Value = 'plot multipal wells on same pad on 3D'
$ws.Range('L10').Value = 'x'
$ws.Range('M10').Value = 'update dirrectional'

$ws.Range('A11').Value = 'X'
$ws.Range('B11').Value = 3.5
$ws.Range('C11').Value = 'Verify north south east and west are labled correctly on 3d plot'
$ws.Range('N11').Value = 'use current import'

$ws.Range('A12').Value = 'X'
$ws.Range('B12').Value = 3.6
$ws.Range('C12').Value = 'annotations on/off button'
$ws.Range('N12').Value = 'if data delete then import new'

$ws.Range('A13').Value = 'X'
$ws.Range('B13').Value = 3.7
$ws.Range('C13').Value = 'Show well check box'
$ws.Range('L13').Value = 'x'
$ws.Range('M13').Value = 'Table view of data?'
$ws.Range('M13').Font.Strikethrough = $true
$ws.Range('N13').Font.Strikethrough = $true

$ws.Range('A14').Value = 'X'
$ws.Range('B14').Value = 3.8
$ws.Range('C14').Value = 'Highlight well check box'
$ws.Range('N14').Value = 'QTableView'
$ws.Range('N14').Font.Strikethrough = $true

$ws.Range('A15').Value = 'X'
$ws.Range('B15').Value = 3.9
$ws.Range('C15').Value = 'Size of chart needs to be bigger 3S-714 example plan view'

# --- Row 16 : "3.1" task (number-formatted) ---
$ws.Range('A16').Value = 'X'
$ws.Range('B16').Value = 3.1
$ws.Range('B16').NumberFormat = "0.00"
$ws.Range('C16').Value = 'uncheck box for formations if none entered'

# --- Row 17 : "3.11" task (struck through / done) ---
$ws.Range('A17').Value = '-'
$ws.Range('B17').Value = 3.11
$ws.Range('B17').Font.Strikethrough = $true
$ws.Range('C17').Value = 'Casing show depths on charts?'
$ws.Range('C17').Font.Strikethrough = $true

# --- Rows 18-25 : renumbered 3.12 - 3.19 ---
$ws.Range('A18').Value = 'X'
$ws.Range('B18').Value = 3.12
$ws.Range('C18').Value = 'Formations names on side of the chart opposite well (if well N/S last point is - put them on the left, if + on right)'

$ws.Range('A19').Value = 'X'
$ws.Range('B19').Value = 3.13
$ws.Range('C19').Value = 'close all charts when closing pad or chart view window'

$ws.Range('A20').Value = 'X'
$ws.Range('B20').Value = 3.14
$ws.Range('C20').Value = 'if show is not selected do not highlight'

$ws.Range('A21').Value = 'X'
$ws.Range('B21').Value = 3.15
$ws.Range('C21').Value = 'annotate wells'

$ws.Range('A22').Value = 'X'
$ws.Range('B22').Value = 3.16
$ws.Range('C22').Value = 'highlight current well'

$ws.Range('A23').Value = 'X'
$ws.Range('B23').Value = 3.17
$ws.Range('C23').Value = 'Add offset well button'

$ws.Range('A24').Value = 'X'
$ws.Range('B24').Value = 3.18
$ws.Range('C24').Value = 'Plan vs Actual window'

$ws.Range('A25').Value = 'X'
$ws.Range('B25').Value = 3.19
$ws.Range('C25').Value = 'Update actual window'

# --- Row 26 : "3.2" task (number-formatted) ---
$ws.Range('A26').Value = 'X'
$ws.Range('B26').Value = 3.2
$ws.Range('B26').NumberFormat = "0.00"
$ws.Range('C26').Value = 'Update plan window'

# --- Rows 27-29 : renumbered 3.21 - 3.23 ---
$ws.Range('A27').Value = 'X'
$ws.Range('B27').Value = 3.21
$ws.Range('C27').Value = 'Modify importCsv.py to set as plan in databse'

$ws.Range('A28').Value = 'X'
$ws.Range('B28').Value = 3.22
$ws.Range('C28').Value = 'Update database to have planned vs actual column'
$ws.Range('D28').Value = 'Under update directional have an add lateral button'

$ws.Range('A29').Value = '-'
$ws.Range('B29').Value = 3.23
$ws.Range('C29').Value = 'Multilateral wells'
$ws.Range('D29').Value = 'Add planed or actual radio putton to latter name window'

# --- Rows 30-34 : 3.23.1 - 3.23.5 sub-tasks ---
$ws.Range('A30').Value = 'X'
$ws.Range('C30').Value = '3.23.1'
$ws.Range('D30').Value = 'Change names in annotations to inclue lateral if not Null'

$ws.Range('A31').Value = 'X'
$ws.Range('C31').Value = '3.23.2'
$ws.Range('D31').Value = 'Chenge colors in chartView.py to black'

$ws.Range('A32').Value = '-'
$ws.Range('C32').Value = '3.23.3'
$ws.Range('D32').Value = 'Make sure planned is not shown in chartView.py if actual exists'

$ws.Range('A33').Value = 'X'
$ws.Range('C33').Value = '3.23.4'

$ws.Range('A34').Value = 'X'
$ws.Range('C34').Value = '3.23.5'

# --- Row 35 : "3.24" TVD vs MD plot task ---
$ws.Range('A35').Value = 'X'
$ws.Range('B35').Value = 3.24
$ws.Range('C35').Value = 'TVD vs MD plot'

# --- Row 36 : new "3.25" task ---
$ws.Range('B36').Value = 3.25
$ws.Range('C36').Value = 'Why does MD vs TVD show a line back to zero when more than one lateral?'

# --- View / selection changes: "2. Dev" becomes the active tab ---
$wellBuilder.Range("D14").Select()
$ws.Activate()
$ws.Range("A18").Select()
